$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.386.33"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.922.62"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0835"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").Value = "3.390.20"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").Value = "2.929.30"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.936"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.29%  "
$ws.Range("D18").Value = "51.391.82"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").Value = "0.0₃0944"
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.70%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  -5.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.52%  "
$ws.Range("E32").Value = "  -7.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("E35").Value = "  -3.43%  "
$ws.Range("E36").Value = "  -5.34%  "
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("E39").Value = "  -8.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.60%  "
$ws.Range("E42").Value = "  -5.58%  "
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.77%  "
$ws.Range("E47").Value = "  +12.12%  "
$ws.Range("D48").Value = "2.017.96"
$ws.Range("E48").Value = "  -4.88%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("D51").Value = "3.212.06"
$ws.Range("E51").Value = "  -2.54%  "
